$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Cells.Item(2, 3).Value = 12595
$ws.Cells.Item(2, 6).Value = 4.5
$ws.Cells.Item(4, 3).Value = 12595
$ws.Cells.Item(4, 6).Value = 4.5
$ws.Cells.Item(6, 3).Value = 12789
$ws.Cells.Item(6, 6).Value = 4.5
$ws.Cells.Item(7, 3).Value = 13303
$ws.Cells.Item(7, 6).Value = 4.99
$ws.Cells.Item(8, 3).Value = 12536
$ws.Cells.Item(8, 6).Value = 4.5
$ws.Cells.Item(9, 3).Value = 12668
$ws.Cells.Item(9, 6).Value = 4.5
$ws.Cells.Item(10, 3).Value = 12527
$ws.Cells.Item(10, 6).Value = 4.5
$ws.Cells.Item(12, 3).Value = 13195
$ws.Cells.Item(12, 6).Value = 4.99
$ws.Cells.Item(14, 3).Value = 12348
$ws.Cells.Item(14, 6).Value = 4.5
$ws.Cells.Item(15, 3).Value = 12500
$ws.Cells.Item(15, 6).Value = 4.5
$ws.Cells.Item(17, 3).Value = 17444
$ws.Cells.Item(17, 6).Value = 4.99
$ws.Cells.Item(18, 3).Value = 13275
$ws.Cells.Item(18, 6).Value = 4.99
$ws.Cells.Item(19, 3).Value = 12613
$ws.Cells.Item(19, 6).Value = 4.5
$ws.Cells.Item(21, 3).Value = 12374
$ws.Cells.Item(21, 6).Value = 4.5
$ws.Cells.Item(23, 3).Value = 13134
$ws.Cells.Item(23, 6).Value = 4.99
$ws.Cells.Item(24, 3).Value = 13354
$ws.Cells.Item(24, 6).Value = 4.99
$ws.Cells.Item(25, 3).Value = 12447
$ws.Cells.Item(25, 6).Value = 4.5
$ws.Cells.Item(27, 3).Value = 12434
$ws.Cells.Item(27, 6).Value = 4.5
$ws.Cells.Item(29, 3).Value = 12389
$ws.Cells.Item(29, 6).Value = 4.5
$ws.Cells.Item(30, 3).Value = 13060
$ws.Cells.Item(30, 6).Value = 4.99
$ws.Cells.Item(35, 3).Value = 13299
$ws.Cells.Item(35, 6).Value = 4.99
$ws.Cells.Item(36, 3).Value = 17448
$ws.Cells.Item(36, 6).Value = 4.99
$ws.Cells.Item(39, 3).Value = 15994
$ws.Cells.Item(39, 6).Value = 4.5
$ws.Cells.Item(40, 3).Value = 12506
$ws.Cells.Item(40, 6).Value = 4.5
$ws.Cells.Item(41, 3).Value = 13406
$ws.Cells.Item(41, 6).Value = 4.99
$ws.Cells.Item(42, 3).Value = 12551
$ws.Cells.Item(42, 6).Value = 4.5
$ws.Cells.Item(43, 3).Value = 12413
$ws.Cells.Item(43, 6).Value = 4.5
$ws.Cells.Item(44, 3).Value = 13422
$ws.Cells.Item(44, 6).Value = 4.99
$ws.Cells.Item(47, 3).Value = 12259
$ws.Cells.Item(47, 6).Value = 4.5
$ws.Cells.Item(51, 3).Value = 12458
$ws.Cells.Item(51, 6).Value = 4.5
$ws.Cells.Item(53, 3).Value = 13330
$ws.Cells.Item(53, 6).Value = 4.99
